$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L170").Value = "Completed"
$ws.Range("M170").Value = "Away Win"
$ws.Range("N170").Value = "Fallo"
$ws.Range("O170").Value = -1.8
$ws.Range("P170").Value = -100.0
$ws.Range("Q170").Value = "2025-09-21 04:26:29"

$ws.Range("L171").Value = "Completed"
$ws.Range("M171").Value = "Home Win"
$ws.Range("N171").Value = "Acierto"
$ws.Range("O171").Value = 1.35
$ws.Range("P171").Value = 45.0
$ws.Range("Q171").Value = "2025-09-21 04:26:29"

$ws.Range("L172").Value = "Completed"
$ws.Range("M172").Value = "Away Win"
$ws.Range("N172").Value = "Fallo"
$ws.Range("O172").Value = -1.7
$ws.Range("P172").Value = -100.0
$ws.Range("Q172").Value = "2025-09-21 04:26:29"

$ws.Range("L173").Value = "Completed"
$ws.Range("M173").Value = "Home Win"
$ws.Range("N173").Value = "Acierto"
$ws.Range("O173").Value = 1.45
$ws.Range("P173").Value = 85.0
$ws.Range("Q173").Value = "2025-09-21 04:26:29"

$ws.Range("L174").Value = "Completed"
$ws.Range("M174").Value = "Home Win"
$ws.Range("N174").Value = "Acierto"
$ws.Range("O174").Value = 1.36
$ws.Range("P174").Value = 105.0
$ws.Range("Q174").Value = "2025-09-21 04:26:29"

$ws.Range("L175").Value = "Completed"
$ws.Range("M175").Value = "Draw"
$ws.Range("N175").Value = "Fallo"
$ws.Range("O175").Value = -0.9
$ws.Range("P175").Value = -100.0
$ws.Range("Q175").Value = "2025-09-21 04:26:29"

$ws.Range("L176").Value = "Completed"
$ws.Range("M176").Value = "Away Win"
$ws.Range("N176").Value = "Fallo"
$ws.Range("O176").Value = -1.0
$ws.Range("P176").Value = -100.0
$ws.Range("Q176").Value = "2025-09-21 04:26:29"

$ws.Range("L177").Value = "Completed"
$ws.Range("M177").Value = "Home Win"
$ws.Range("N177").Value = "Acierto"
$ws.Range("O177").Value = 1.47
$ws.Range("P177").Value = 70.0
$ws.Range("Q177").Value = "2025-09-21 04:26:29"

$ws.Range("L178").Value = "Completed"
$ws.Range("M178").Value = "Draw"
$ws.Range("N178").Value = "Fallo"
$ws.Range("O178").Value = -1.1
$ws.Range("P178").Value = -100.0
$ws.Range("Q178").Value = "2025-09-21 04:26:29"

$ws.Range("L179").Value = "Completed"
$ws.Range("M179").Value = "Home Win"
$ws.Range("N179").Value = "Acierto"
$ws.Range("O179").Value = 1.23
$ws.Range("P179").Value = 95.0
$ws.Range("Q179").Value = "2025-09-21 04:26:29"

$ws.Range("L180").Value = "Completed"
$ws.Range("M180").Value = "Draw"
$ws.Range("N180").Value = "Fallo"
$ws.Range("O180").Value = -3.0
$ws.Range("P180").Value = -100.0
$ws.Range("Q180").Value = "2025-09-21 04:26:29"

$ws.Range("L181").Value = "Completed"
$ws.Range("M181").Value = "Home Win"
$ws.Range("N181").Value = "Acierto"
$ws.Range("O181").Value = 1.88
$ws.Range("P181").Value = 65.0
$ws.Range("Q181").Value = "2025-09-21 04:26:29"

$ws.Range("L182").Value = "Completed"
$ws.Range("M182").Value = "Home Win"
$ws.Range("N182").Value = "Acierto"
$ws.Range("O182").Value = 1.59
$ws.Range("P182").Value = 53.0
$ws.Range("Q182").Value = "2025-09-21 04:26:29"

$ws.Range("L183").Value = "Completed"
$ws.Range("M183").Value = "Home Win"
$ws.Range("N183").Value = "Acierto"
$ws.Range("O183").Value = 1.53
$ws.Range("P183").Value = 73.0
$ws.Range("Q183").Value = "2025-09-21 04:26:29"

$ws.Range("L184").Value = "Completed"
$ws.Range("M184").Value = "Draw"
$ws.Range("N184").Value = "Fallo"
$ws.Range("O184").Value = -0.6
$ws.Range("P184").Value = -100.0
$ws.Range("Q184").Value = "2025-09-21 04:26:29"

$ws.Range("L185").Value = "Completed"
$ws.Range("M185").Value = "Draw"
$ws.Range("N185").Value = "Fallo"
$ws.Range("O185").Value = -2.1
$ws.Range("P185").Value = -100.0
$ws.Range("Q185").Value = "2025-09-21 04:26:29"

$ws.Range("L186").Value = "Completed"
$ws.Range("M186").Value = "Home Win"
$ws.Range("N186").Value = "Fallo"
$ws.Range("O186").Value = -1.7
$ws.Range("P186").Value = -100.0
$ws.Range("Q186").Value = "2025-09-21 04:26:29"

$ws.Range("L187").Value = "Completed"
$ws.Range("M187").Value = "Home Win"
$ws.Range("N187").Value = "Acierto"
$ws.Range("O187").Value = 1.59
$ws.Range("P187").Value = 53.0
$ws.Range("Q187").Value = "2025-09-21 04:26:29"

$ws.Range("L188").Value = "Completed"
$ws.Range("M188").Value = "Away Win"
$ws.Range("N188").Value = "Acierto"
$ws.Range("O188").Value = 1.74
$ws.Range("P188").Value = 83.0
$ws.Range("Q188").Value = "2025-09-21 04:26:29"

$ws.Range("L189").Value = "Completed"
$ws.Range("M189").Value = "Draw"
$ws.Range("N189").Value = "Fallo"
$ws.Range("O189").Value = -2.2
$ws.Range("P189").Value = -100.0
$ws.Range("Q189").Value = "2025-09-21 04:26:29"

$ws.Range("L190").Value = "Completed"
$ws.Range("M190").Value = "Home Win"
$ws.Range("N190").Value = "Acierto"
$ws.Range("O190").Value = 1.5
$ws.Range("P190").Value = 100.0
$ws.Range("Q190").Value = "2025-09-21 15:20:12"

$ws.Range("L191").Value = "Completed"
$ws.Range("M191").Value = "Home Win"
$ws.Range("N191").Value = "Acierto"
$ws.Range("O191").Value = 1.21
$ws.Range("P191").Value = 110.0
$ws.Range("Q191").Value = "2025-09-21 15:20:12"

$ws.Range("L192").Value = "Completed"
$ws.Range("M192").Value = "Home Win"
$ws.Range("N192").Value = "Fallo"
$ws.Range("O192").Value = -0.3
$ws.Range("P192").Value = -100.0
$ws.Range("Q192").Value = "2025-09-21 15:20:12"

$ws.Range("L193").Value = "Completed"
$ws.Range("M193").Value = "Draw"
$ws.Range("N193").Value = "Fallo"
$ws.Range("O193").Value = -0.8
$ws.Range("P193").Value = -100.0
$ws.Range("Q193").Value = "2025-09-21 15:20:12"

$ws.Range("L194").Value = "Completed"
$ws.Range("M194").Value = "Draw"
$ws.Range("N194").Value = "Fallo"
$ws.Range("O194").Value = -3.0
$ws.Range("P194").Value = -100.0
$ws.Range("Q194").Value = "2025-09-21 15:20:12"

$ws.Range("L195").Value = "Completed"
$ws.Range("M195").Value = "Away Win"
$ws.Range("N195").Value = "Acierto"
$ws.Range("O195").Value = 0.52
$ws.Range("P195").Value = 130.0
$ws.Range("Q195").Value = "2025-09-21 15:20:12"

$ws.Range("L196").Value = "Completed"
$ws.Range("M196").Value = "Home Win"
$ws.Range("N196").Value = "Acierto"
$ws.Range("O196").Value = 1.1
$ws.Range("P196").Value = 110.0
$ws.Range("Q196").Value = "2025-09-21 15:20:12"
